# NetLiquidity FRED data refresh: append new WTREGEN observation and
# update the SeriesInfo metadata block to match.

$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("Data")
$infoWs = $wb.Worksheets.Item("SeriesInfo")

# --- Data sheet: append the new weekly observation (row 94) ---
$dataWs.Range("A93").Copy($dataWs.Range("A94"))
$dataWs.Range("A94").Value = 45119
$dataWs.Range("B94").Value = 514.337

# --- SeriesInfo sheet: refresh the FRED metadata timestamps ---
# Force the plain YYYY-MM-DD values to stay text (Excel would otherwise
# auto-convert them to date serials), then drop the text number format so
# the cell keeps the workbook's default (unstyled) appearance.
$infoWs.Range("B3").NumberFormat = "@"
$infoWs.Range("B3").Value = "2023-07-20"
$infoWs.Range("B3").ClearFormats()

$infoWs.Range("B4").NumberFormat = "@"
$infoWs.Range("B4").Value = "2023-07-20"
$infoWs.Range("B4").ClearFormats()

$infoWs.Range("B7").NumberFormat = "@"
$infoWs.Range("B7").Value = "2023-07-12"
$infoWs.Range("B7").ClearFormats()

$infoWs.Range("B14").Value = "2023-07-13 15:35:18-05"
